$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.423.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.648.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.98%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.640.87"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.612"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.50%  "
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.198"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.83"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +24.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.604"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "48.59"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000286"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.230.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "670.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.88"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.645.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.474.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.931"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "100.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.88%  "
$ws.Range("E31").Value = "  -2.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.60%  "
$ws.Range("E33").Value = "  -4.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.32"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.96"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "580.52"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.05"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.107"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "58.06"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.582.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.83%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0453"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.41%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.141"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.344"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "34.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₃0743"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("E47").Value = "  -1.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.94"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.133"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "135.46"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.57%  "
